# Apply weekly update: insert two new data rows (Primera / Segunda quality
# records for Apio, Americana, Terminal La Palmera de La Serena) at the top
# of the data block (just before the former row 417), shifting the
# existing rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 417.
$ws.Rows("417:418").Insert()

# Row 417: new "Primera" quality record
$ws.Cells.Item(417, 1).Value2 = 8
$ws.Cells.Item(417, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(417, 3).Value2 = "Coquimbo"
$ws.Cells.Item(417, 4).Value2 = 44769
$ws.Cells.Item(417, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(417, 5).Value2 = 4
$ws.Cells.Item(417, 6).Value2 = 100112017
$ws.Cells.Item(417, 7).Value2 = "Apio"
$ws.Cells.Item(417, 8).Value2 = "Americana (o)"
$ws.Cells.Item(417, 9).Value2 = "Primera"
$ws.Cells.Item(417, 10).Value2 = 1600
$ws.Cells.Item(417, 11).Value2 = 8000
$ws.Cells.Item(417, 12).Value2 = 9000
$ws.Cells.Item(417, 13).Value2 = 8500
$ws.Cells.Item(417, 14).Value2 = "$/docena de matas"
$ws.Cells.Item(417, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(417, 16).Value2 = 1417
$ws.Cells.Item(417, 17).Value2 = 6
$ws.Cells.Item(417, 18).Value2 = "Hortaliza"

# Row 418: new "Segunda" quality record
$ws.Cells.Item(418, 1).Value2 = 8
$ws.Cells.Item(418, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(418, 3).Value2 = "Coquimbo"
$ws.Cells.Item(418, 4).Value2 = 44769
$ws.Cells.Item(418, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(418, 5).Value2 = 4
$ws.Cells.Item(418, 6).Value2 = 100112017
$ws.Cells.Item(418, 7).Value2 = "Apio"
$ws.Cells.Item(418, 8).Value2 = "Americana (o)"
$ws.Cells.Item(418, 9).Value2 = "Segunda"
$ws.Cells.Item(418, 10).Value2 = 1200
$ws.Cells.Item(418, 11).Value2 = 6500
$ws.Cells.Item(418, 12).Value2 = 7000
$ws.Cells.Item(418, 13).Value2 = 6750
$ws.Cells.Item(418, 14).Value2 = "$/docena de matas"
$ws.Cells.Item(418, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(418, 16).Value2 = 1125
$ws.Cells.Item(418, 17).Value2 = 6
$ws.Cells.Item(418, 18).Value2 = "Hortaliza"
